$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new "dataset" table appended below the existing content.
$ws.Cells.Item(145, 1).Value = "数据集"
$ws.Cells.Item(145, 2).Value = "域"
$ws.Cells.Item(145, 3).Value = "N"
$ws.Cells.Item(145, 4).Value = "M"
$ws.Cells.Item(145, 5).Value = "L "
$ws.Cells.Item(145, 6).Value = "LC(D)"
$ws.Cells.Item(145, 7).Value = "LD(D)"

# Data rows 146-158: dataset name, domain, N, M, L, LC(D), LD(D)
$ws.Cells.Item(146, 1).Value = "MediaMill"
$ws.Cells.Item(146, 2).Value = "视频"
$ws.Cells.Item(146, 3).Value = 43907
$ws.Cells.Item(146, 4).Value = 120
$ws.Cells.Item(146, 5).Value = 101
$ws.Cells.Item(146, 6).Value = 4.376
$ws.Cells.Item(146, 7).Value = 0.044

$ws.Cells.Item(147, 1).Value = "TMC2007"
$ws.Cells.Item(147, 2).Value = "文本"
$ws.Cells.Item(147, 3).Value = 28596
$ws.Cells.Item(147, 4).Value = 500
$ws.Cells.Item(147, 5).Value = 22
$ws.Cells.Item(147, 6).Value = 2.16
$ws.Cells.Item(147, 7).Value = 0.098

$ws.Cells.Item(148, 1).Value = "Rcv1-v2"
$ws.Cells.Item(148, 2).Value = "文本"
$ws.Cells.Item(148, 3).Value = 804414
$ws.Cells.Item(148, 4).Value = 500
$ws.Cells.Item(148, 5).Value = 103
$ws.Cells.Item(148, 6).Value = 3.24
$ws.Cells.Item(148, 7).Value = 0.031

$ws.Cells.Item(149, 1).Value = "IMDB"
$ws.Cells.Item(149, 2).Value = "文本"
$ws.Cells.Item(149, 3).Value = 120919
$ws.Cells.Item(149, 4).Value = 1001
$ws.Cells.Item(149, 5).Value = 28
$ws.Cells.Item(149, 6).Value = 2
$ws.Cells.Item(149, 7).Value = 0.071

$ws.Cells.Item(150, 1).Value = "20NG"
$ws.Cells.Item(150, 2).Value = "文本"
$ws.Cells.Item(150, 3).Value = 19300
$ws.Cells.Item(150, 4).Value = 1006
$ws.Cells.Item(150, 5).Value = 20
$ws.Cells.Item(150, 6).Value = 1.02
$ws.Cells.Item(150, 7).Value = 0.051

$ws.Cells.Item(151, 1).Value = "Yeast"
$ws.Cells.Item(151, 2).Value = "生物"
$ws.Cells.Item(151, 3).Value = 2417
$ws.Cells.Item(151, 4).Value = 103
$ws.Cells.Item(151, 5).Value = 14
$ws.Cells.Item(151, 6).Value = 4.237
$ws.Cells.Item(151, 7).Value = 0.303

$ws.Cells.Item(152, 1).Value = "Ohsumed"
$ws.Cells.Item(152, 2).Value = "文本"
$ws.Cells.Item(152, 3).Value = 13529
$ws.Cells.Item(152, 4).Value = 1002
$ws.Cells.Item(152, 5).Value = 23
$ws.Cells.Item(152, 6).Value = 1.66
$ws.Cells.Item(152, 7).Value = 0.072

$ws.Cells.Item(153, 1).Value = "Slashdot"
$ws.Cells.Item(153, 2).Value = "文本"
$ws.Cells.Item(153, 3).Value = 3782
$ws.Cells.Item(153, 4).Value = 1079
$ws.Cells.Item(153, 5).Value = 22
$ws.Cells.Item(153, 6).Value = 1.18
$ws.Cells.Item(153, 7).Value = 0.053

$ws.Cells.Item(154, 1).Value = "Reuters"
$ws.Cells.Item(154, 2).Value = "文本"
$ws.Cells.Item(154, 3).Value = 6000
$ws.Cells.Item(154, 4).Value = 500
$ws.Cells.Item(154, 5).Value = 101
$ws.Cells.Item(154, 6).Value = 2.88
$ws.Cells.Item(154, 7).Value = 0.028

$ws.Cells.Item(155, 1).Value = "Enron"
$ws.Cells.Item(155, 2).Value = "文本"
$ws.Cells.Item(155, 3).Value = 1702
$ws.Cells.Item(155, 4).Value = 1001
$ws.Cells.Item(155, 5).Value = 53
$ws.Cells.Item(155, 6).Value = 3.4
$ws.Cells.Item(155, 7).Value = 0.064

$ws.Cells.Item(156, 1).Value = "Scene"
$ws.Cells.Item(156, 2).Value = "图片"
$ws.Cells.Item(156, 3).Value = 2407
$ws.Cells.Item(156, 4).Value = 294
$ws.Cells.Item(156, 5).Value = 6
$ws.Cells.Item(156, 6).Value = 1.074
$ws.Cells.Item(156, 7).Value = 0.179

$ws.Cells.Item(157, 1).Value = "Medical"
$ws.Cells.Item(157, 2).Value = "文本"
$ws.Cells.Item(157, 3).Value = 978
$ws.Cells.Item(157, 4).Value = 1449
$ws.Cells.Item(157, 5).Value = 45
$ws.Cells.Item(157, 6).Value = 1.25
$ws.Cells.Item(157, 7).Value = 0.028

$ws.Cells.Item(158, 1).Value = "Core15K"
$ws.Cells.Item(158, 2).Value = "多媒体"
$ws.Cells.Item(158, 3).Value = 5000
$ws.Cells.Item(158, 4).Value = 499
$ws.Cells.Item(158, 5).Value = 374
$ws.Cells.Item(158, 6).Value = 3.52
$ws.Cells.Item(158, 7).Value = 0.009

# Match the saved view state: scrolled down with C146 the active cell.
$ws.Range("C146").Select()
